$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("C2").Value = 1.629813710267143
$ws.Range("E2").Value = 3.259627420534287
$ws.Range("F2").Value = 2.04009049

# Row 3
$ws.Range("C3").Value = 1.381725987173144
$ws.Range("E3").Value = 5.526903948692576
$ws.Range("F3").Value = 6.387462731273127

# Row 4
$ws.Range("C4").Value = 0
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = 0
$ws.Range("F4").Value = 0

# Row 5
$ws.Range("C5").Value = 1.443728898093037
$ws.Range("E5").Value = 8.662373388558221
$ws.Range("F5").Value = 4.313260548277049

# Row 6
$ws.Range("C6").Value = 1.390120801212101
$ws.Range("D6").Value = 3.183609109940719
$ws.Range("E6").Value = 8.340724807272604
$ws.Range("F6").Value = 5.815722961037232

# Row 7
$ws.Range("C7").Value = 0
$ws.Range("D7").Value = 0
$ws.Range("E7").Value = 0
$ws.Range("F7").Value = 0

# Row 8
$ws.Range("C8").Value = 0
$ws.Range("D8").Value = 0
$ws.Range("E8").Value = 0
$ws.Range("F8").Value = 0
$ws.Range("G8").Value = "lipid/free"

# Row 9
$ws.Range("C9").Value = 0
$ws.Range("D9").Value = 0
$ws.Range("E9").Value = 0
$ws.Range("F9").Value = 0

# Row 10
$ws.Range("C10").Value = 0.8335219886725862
$ws.Range("E10").Value = 3.334087954690345
$ws.Range("F10").Value = 0.952300540779256

# Row 11
$ws.Range("C11").Value = 0
$ws.Range("D11").Value = 0
$ws.Range("E11").Value = 0
$ws.Range("F11").Value = 0

# Row 12
$ws.Range("C12").Value = 0.7363978474716336
$ws.Range("D12").Value = 2.945591389886534
$ws.Range("E12").Value = 2.945591389886534
$ws.Range("F12").Value = 0.1432853375000008

# Row 13
$ws.Range("C13").Value = 0
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("F13").Value = 0
